$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 160
$ws.Range("I6").Value = 160
$ws.Range("K6").Value = 480
$ws.Range("M6").Value = -368
$ws.Range("H12").Value = 1416.3334
$ws.Range("I12").Value = 1249.5
$ws.Range("K12").Value = 1249.5
$ws.Range("M12").Value = -1079.5
$ws.Range("H29").Value = 1919.6
$ws.Range("I29").Value = 499
$ws.Range("J29").Value = 2866.6667
$ws.Range("K29").Value = 1497
$ws.Range("L29").Value = 8600.000100000001
$ws.Range("N29").Value = -9162.000100000001
$ws.Range("M29").Value = -1216
$ws.Range("H38").Value = 2898.6316
$ws.Range("I38").Value = 218.33333
$ws.Range("J38").Value = 4135.6924
$ws.Range("K38").Value = 654.99999
$ws.Range("L38").Value = 12407.0772
$ws.Range("M38").Value = -282.99999
$ws.Range("N38").Value = -13151.0772
$ws.Range("H40").Value = 4080.3333
$ws.Range("I40").Value = 2349.9167
$ws.Range("J40").Value = 5069.143
$ws.Range("K40").Value = 2349.9167
$ws.Range("L40").Value = 5069.143
$ws.Range("M40").Value = -2174.9167
$ws.Range("N40").Value = -5419.143
$ws.Range("H58").Value = 149.8
$ws.Range("I58").Value = 62.25
$ws.Range("J58").Value = 500
$ws.Range("K58").Value = 186.75
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -36.75
$ws.Range("N58").Value = -1800
$ws.Range("H87").Value = 68313
$ws.Range("J87").Value = 68313
$ws.Range("L87").Value = 68313
$ws.Range("N87").Value = -70809
$ws.Range("H90").Value = 68313
$ws.Range("J90").Value = 68313
$ws.Range("L90").Value = 204939
$ws.Range("N90").Value = -217419
$ws.Range("H92").Value = 734.1
$ws.Range("I92").Value = 787.8889
$ws.Range("K92").Value = 787.8889
$ws.Range("M92").Value = 460.1111
$ws.Range("H97").Value = 2195.7144
$ws.Range("J97").Value = 2195.7144
$ws.Range("L97").Value = 6587.1432
$ws.Range("N97").Value = -7579.1432
$ws.Range("H98").Value = 2099.0715
$ws.Range("I98").Value = 1035.125
$ws.Range("J98").Value = 3517.6667
$ws.Range("K98").Value = 1035.125
$ws.Range("L98").Value = 3517.6667
$ws.Range("M98").Value = 462.875
$ws.Range("N98").Value = -6513.6667
$ws.Range("H100").Value = 1887
$ws.Range("I100").Value = 1941.625
$ws.Range("J100").Value = 1450
$ws.Range("K100").Value = 1941.625
$ws.Range("L100").Value = 1450
$ws.Range("M100").Value = -1400.625
$ws.Range("N100").Value = -2532
$ws.Range("H107").Value = 457.6
$ws.Range("I107").Value = 420.35294
$ws.Range("K107").Value = 420.35294
$ws.Range("M107").Value = 1499.64706
$ws.Range("H113").Value = 6987.75
$ws.Range("J113").Value = 6987.75
$ws.Range("L113").Value = 6987.75
$ws.Range("N113").Value = -13495.75
$ws.Range("H116").Value = 7175.3887
$ws.Range("I116").Value = 6515.636
$ws.Range("K116").Value = 6515.636
$ws.Range("M116").Value = -3073.636
$ws.Range("H122").Value = 2099.0715
$ws.Range("I122").Value = 1035.125
$ws.Range("J122").Value = 3517.6667
$ws.Range("K122").Value = 3105.375
$ws.Range("L122").Value = 10553.0001
$ws.Range("M122").Value = -655.375
$ws.Range("N122").Value = -15453.0001
$ws.Range("H132").Value = 3329.2727
$ws.Range("I132").Value = 3329.2727
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9987.8181
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7457.8181
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 3477.9583
$ws.Range("J137").Value = 3958.5293
$ws.Range("L137").Value = 11875.5879
$ws.Range("N137").Value = -16975.5879
$ws.Range("H138").Value = 5096.25
$ws.Range("I138").Value = 3495.5881
$ws.Range("J138").Value = 6910.3335
$ws.Range("K138").Value = 10486.7643
$ws.Range("L138").Value = 20731.0005
$ws.Range("M138").Value = -5346.764299999999
$ws.Range("N138").Value = -31011.0005
$ws.Range("H140").Value = 94398
$ws.Range("J140").Value = 94398
$ws.Range("L140").Value = 94398
$ws.Range("N140").Value = -104758
$ws.Range("H141").Value = 4712.154
$ws.Range("I141").Value = 3731.75
$ws.Range("J141").Value = 6280.8
$ws.Range("K141").Value = 11195.25
$ws.Range("L141").Value = 18842.4
$ws.Range("M141").Value = -6015.25
$ws.Range("N141").Value = -29202.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 246
$ws.Range("I5").Value = 361.5
$ws.Range("K5").Value = 361.5
$ws.Range("M5").Value = -249.5
$ws.Range("H32").Value = 6430.9834
$ws.Range("I32").Value = 2585.8518
$ws.Range("K32").Value = 2585.8518
$ws.Range("M32").Value = -2298.8518
$ws.Range("H45").Value = 1574.875
$ws.Range("I45").Value = 1479.8667
$ws.Range("K45").Value = 1479.8667
$ws.Range("M45").Value = -1102.8667
$ws.Range("H61").Value = 4973.35
$ws.Range("I61").Value = 1905.5
$ws.Range("K61").Value = 1905.5
$ws.Range("M61").Value = -1693.5
$ws.Range("H74").Value = 8407.833000000001
$ws.Range("I74").Value = 3698
$ws.Range("J74").Value = 10762.75
$ws.Range("K74").Value = 3698
$ws.Range("L74").Value = 10762.75
$ws.Range("M74").Value = -2824
$ws.Range("N74").Value = -12510.75
$ws.Range("H77").Value = 8407.833000000001
$ws.Range("I77").Value = 3698
$ws.Range("J77").Value = 10762.75
$ws.Range("K77").Value = 18490
$ws.Range("L77").Value = 53813.75
$ws.Range("M77").Value = -14122
$ws.Range("N77").Value = -62549.75
$ws.Range("H110").Value = 2248.1177
$ws.Range("I110").Value = 2307.375
$ws.Range("K110").Value = 2307.375
$ws.Range("M110").Value = -262.375
$ws.Range("H122").Value = 4596.75
$ws.Range("I122").Value = 4572.7393
$ws.Range("J122").Value = 5149
$ws.Range("K122").Value = 13718.2179
$ws.Range("L122").Value = 15447
$ws.Range("M122").Value = -11268.2179
$ws.Range("N122").Value = -20347
$ws.Range("H132").Value = 3357.8965
$ws.Range("I132").Value = 2493.8635
$ws.Range("J132").Value = 6073.4287
$ws.Range("K132").Value = 7481.5905
$ws.Range("L132").Value = 18220.2861
$ws.Range("M132").Value = -4951.5905
$ws.Range("N132").Value = -23280.2861
$ws.Range("H136").Value = 4973.35
$ws.Range("I136").Value = 1905.5
$ws.Range("K136").Value = 5716.5
$ws.Range("M136").Value = -3166.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 246
$ws.Range("I4").Value = 361.5
$ws.Range("K4").Value = 361.5
$ws.Range("M4").Value = -246.5
$ws.Range("H20").Value = 1448.5172
$ws.Range("I20").Value = 1581.25
$ws.Range("J20").Value = 811.4
$ws.Range("K20").Value = 1581.25
$ws.Range("L20").Value = 811.4
$ws.Range("M20").Value = -1334.25
$ws.Range("N20").Value = -1305.4
$ws.Range("H22").Value = 672
$ws.Range("I22").Value = 589.6667
$ws.Range("K22").Value = 589.6667
$ws.Range("M22").Value = -416.6667
$ws.Range("H80").Value = 48032.777
$ws.Range("J80").Value = 18899.143
$ws.Range("L80").Value = 18899.143
$ws.Range("N80").Value = -20895.143
$ws.Range("H83").Value = 48032.777
$ws.Range("J83").Value = 18899.143
$ws.Range("L83").Value = 94495.715
$ws.Range("N83").Value = -104479.715
$ws.Range("H86").Value = 3021.875
$ws.Range("I86").Value = 1838.2
$ws.Range("K86").Value = 1838.2
$ws.Range("M86").Value = -715.2
$ws.Range("H89").Value = 3021.875
$ws.Range("I89").Value = 1838.2
$ws.Range("K89").Value = 9191
$ws.Range("M89").Value = -3575
$ws.Range("H107").Value = 2476.0833
$ws.Range("J107").Value = 3266
$ws.Range("L107").Value = 3266
$ws.Range("N107").Value = -7106
$ws.Range("H134").Value = 3348.9167
$ws.Range("I134").Value = 2055.4443
$ws.Range("K134").Value = 6166.3329
$ws.Range("M134").Value = -3631.3329
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8497.652
$ws.Range("I31").Value = 4125.5
$ws.Range("K31").Value = 4125.5
$ws.Range("M31").Value = -3830.5
$ws.Range("H34").Value = 8497.652
$ws.Range("I34").Value = 4125.5
$ws.Range("K34").Value = 4125.5
$ws.Range("M34").Value = -3923.5
$ws.Range("H41").Value = 24592.445
$ws.Range("I41").Value = 9867.799999999999
$ws.Range("J41").Value = 42998.25
$ws.Range("K41").Value = 9867.799999999999
$ws.Range("L41").Value = 42998.25
$ws.Range("M41").Value = -9439.799999999999
$ws.Range("N41").Value = -43854.25
$ws.Range("H60").Value = 22122.154
$ws.Range("I60").Value = 14899.667
$ws.Range("J60").Value = 28312.857
$ws.Range("K60").Value = 14899.667
$ws.Range("L60").Value = 28312.857
$ws.Range("M60").Value = -14388.667
$ws.Range("N60").Value = -29334.857
$ws.Range("H99").Value = 2499.6428
$ws.Range("I99").Value = 1999.6154
$ws.Range("J99").Value = 9000
$ws.Range("K99").Value = 1999.6154
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = -501.6153999999999
$ws.Range("N99").Value = -11996
$ws.Range("H126").Value = 2499.6428
$ws.Range("I126").Value = 1999.6154
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 5998.8462
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -3528.8462
$ws.Range("N126").Value = -31940
$ws.Range("H132").Value = 3671.7144
$ws.Range("I132").Value = 2661.4285
$ws.Range("K132").Value = 7984.2855
$ws.Range("M132").Value = -5454.2855
$ws.Range("H134").Value = 5208.3887
$ws.Range("I134").Value = 3068.5
$ws.Range("J134").Value = 12698
$ws.Range("K134").Value = 9205.5
$ws.Range("L134").Value = 38094
$ws.Range("M134").Value = -6670.5
$ws.Range("N134").Value = -43164
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2884.2778
$ws.Range("I3").Value = 2739.625
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 8218.875
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = -8106.875
$ws.Range("N3").Value = -9224
$ws.Range("H4").Value = 26043438
$ws.Range("I4").Value = 104167620
$ws.Range("J4").Value = 2044.5
$ws.Range("K4").Value = 312502860
$ws.Range("L4").Value = 6133.5
$ws.Range("M4").Value = -312502748
$ws.Range("N4").Value = -6357.5
$ws.Range("H107").Value = 572.6667
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H123").Value = 8432.25
$ws.Range("I123").Value = 8432.25
$ws.Range("K123").Value = 25296.75
$ws.Range("M123").Value = -22846.75
$ws.Range("H140").Value = 1639.9025
$ws.Range("J140").Value = 1635.4517
$ws.Range("L140").Value = 4906.355100000001
$ws.Range("N140").Value = -15266.3551
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 12403.134
$ws.Range("I3").Value = 19835.428
$ws.Range("J3").Value = 5899.875
$ws.Range("K3").Value = 19835.428
$ws.Range("L3").Value = 5899.875
$ws.Range("M3").Value = -19719.428
$ws.Range("N3").Value = -6131.875
$ws.Range("H46").Value = 35015.5
$ws.Range("J46").Value = 69990
$ws.Range("L46").Value = 69990
$ws.Range("N46").Value = -70302
$ws.Range("H57").Value = 54500
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 54500
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 54500
$ws.Range("N57").Value = -56140
$ws.Range("M57").ClearContents()
$ws.Range("H70").Value = 4271.8613
$ws.Range("I70").Value = 2531.1
$ws.Range("J70").Value = 4941.385
$ws.Range("K70").Value = 2531.1
$ws.Range("L70").Value = 4941.385
$ws.Range("M70").Value = -2261.1
$ws.Range("N70").Value = -5481.385
$ws.Range("H73").Value = 4271.8613
$ws.Range("I73").Value = 2531.1
$ws.Range("J73").Value = 4941.385
$ws.Range("K73").Value = 2531.1
$ws.Range("L73").Value = 4941.385
$ws.Range("M73").Value = -1595.1
$ws.Range("N73").Value = -6813.385
$ws.Range("H80").Value = 12957.8
$ws.Range("I80").Value = 12447.25
$ws.Range("J80").Value = 15000
$ws.Range("K80").Value = 12447.25
$ws.Range("L80").Value = 15000
$ws.Range("M80").Value = -11449.25
$ws.Range("N80").Value = -16996
$ws.Range("H83").Value = 12957.8
$ws.Range("I83").Value = 12447.25
$ws.Range("J83").Value = 15000
$ws.Range("K83").Value = 62236.25
$ws.Range("L83").Value = 75000
$ws.Range("M83").Value = -57244.25
$ws.Range("N83").Value = -84984
$ws.Range("H122").Value = 6265.923
$ws.Range("I122").Value = 4430
$ws.Range("J122").Value = 13976.8
$ws.Range("K122").Value = 13290
$ws.Range("L122").Value = 41930.39999999999
$ws.Range("M122").Value = -10840
$ws.Range("N122").Value = -46830.39999999999
$ws.Range("H126").Value = 4748.6665
$ws.Range("I126").Value = 2998.2222
$ws.Range("K126").Value = 8994.6666
$ws.Range("M126").Value = -6524.6666
$ws.Range("H132").Value = 7018.171
$ws.Range("I132").Value = 6861.7026
$ws.Range("J132").Value = 8465.5
$ws.Range("K132").Value = 20585.1078
$ws.Range("L132").Value = 25396.5
$ws.Range("M132").Value = -18055.1078
$ws.Range("N132").Value = -30456.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4468.273
$ws.Range("I40").Value = 3115.1
$ws.Range("K40").Value = 3115.1
$ws.Range("M40").Value = -2979.1
$ws.Range("H46").Value = 4615
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 4782.727
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 4782.727
$ws.Range("M46").Value = -3812
$ws.Range("N46").Value = -5158.727
$ws.Range("H55").Value = 1786.3334
$ws.Range("I55").Value = 495
$ws.Range("K55").Value = 495
$ws.Range("M55").Value = -322
$ws.Range("H82").Value = 1227.2778
$ws.Range("I82").Value = 856.8570999999999
$ws.Range("J82").Value = 1463
$ws.Range("K82").Value = 856.8570999999999
$ws.Range("L82").Value = 1463
$ws.Range("M82").Value = -495.8570999999999
$ws.Range("N82").Value = -2185
$ws.Range("H85").Value = 1227.2778
$ws.Range("I85").Value = 856.8570999999999
$ws.Range("J85").Value = 1463
$ws.Range("K85").Value = 856.8570999999999
$ws.Range("L85").Value = 1463
$ws.Range("M85").Value = 391.1429000000001
$ws.Range("N85").Value = -3959
$ws.Range("H122").Value = 6285
$ws.Range("I122").Value = 4999.25
$ws.Range("J122").Value = 13999.5
$ws.Range("K122").Value = 14997.75
$ws.Range("L122").Value = 41998.5
$ws.Range("M122").Value = -12547.75
$ws.Range("N122").Value = -46898.5
$ws.Range("H132").Value = 7004.9062
$ws.Range("I132").Value = 6072.087
$ws.Range("J132").Value = 9388.777
$ws.Range("K132").Value = 18216.261
$ws.Range("L132").Value = 28166.331
$ws.Range("M132").Value = -15686.261
$ws.Range("N132").Value = -33226.331
$ws.Range("H136").Value = 6686.8423
$ws.Range("I136").Value = 5345.2
$ws.Range("J136").Value = 8177.5557
$ws.Range("K136").Value = 16035.6
$ws.Range("L136").Value = 24532.6671
$ws.Range("M136").Value = -13485.6
$ws.Range("N136").Value = -29632.6671
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 91665.5
$ws.Range("J63").Value = 99998.60000000001
$ws.Range("L63").Value = 99998.60000000001
$ws.Range("N63").Value = -101246.6
$ws.Range("H66").Value = 91665.5
$ws.Range("J66").Value = 99998.60000000001
$ws.Range("L66").Value = 299995.8
$ws.Range("N66").Value = -306235.8
$ws.Range("H122").Value = 3798.9092
$ws.Range("J122").Value = 3392.3333
$ws.Range("L122").Value = 10176.9999
$ws.Range("N122").Value = -15076.9999
$ws.Range("H132").Value = 3528.9375
$ws.Range("I132").Value = 2113.25
$ws.Range("J132").Value = 4944.625
$ws.Range("K132").Value = 6339.75
$ws.Range("L132").Value = 14833.875
$ws.Range("M132").Value = -3809.75
$ws.Range("N132").Value = -19893.875
$ws.Range("H136").Value = 7579.591
$ws.Range("I136").Value = 7646.6665
$ws.Range("J136").Value = 7435.857
$ws.Range("K136").Value = 22939.9995
$ws.Range("L136").Value = 22307.571
$ws.Range("M136").Value = -20389.9995
$ws.Range("N136").Value = -27407.571
